$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.060.28"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.819.95"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "'233.21"
$ws.Range("E5").Value = "  -1.64%  "

$ws.Range("D6").Value = "'0.5911"
$ws.Range("E6").Value = "  -3.04%  "

$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").Value = "'0.2743"
$ws.Range("E8").Value = "  -2.77%  "

$ws.Range("D9").Value = "'0.06786"
$ws.Range("E9").Value = "  -4.35%  "

$ws.Range("D10").Value = "'22.94"
$ws.Range("E10").Value = "  -4.06%  "

$ws.Range("D11").Value = "'0.07505"
$ws.Range("E11").Value = "  -1.87%  "

$ws.Range("D12").Value = "1.816.12"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").Value = "'4.673"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").Value = "'0.6229"
$ws.Range("E14").Value = "  -1.96%  "

$ws.Range("D15").Value = "'0.000009407"
$ws.Range("E15").Value = "  -6.00%  "

$ws.Range("D16").Value = "'74.52"
$ws.Range("E16").Value = "  -6.44%  "

$ws.Range("D17").Value = "28.791.83"
$ws.Range("E17").Value = "  -1.44%  "

$ws.Range("D18").Value = "'5.422"
$ws.Range("E18").Value = "  -9.07%  "

$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "'207.73"
$ws.Range("E20").Value = "  -9.38%  "

$ws.Range("D21").Value = "'11.37"
$ws.Range("E21").Value = "  -3.81%  "

$ws.Range("D22").Value = "'6.767"
$ws.Range("E22").Value = "  -3.81%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'154.14"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("D25").Value = "'0.1269"
$ws.Range("E25").Value = "  -2.30%  "

$ws.Range("E26").Value = "  -4.03%  "

$ws.Range("D27").Value = "'16.25"
$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("D28").Value = "'0.06406"
$ws.Range("E28").Value = "  -5.71%  "

$ws.Range("D29").Value = "'1.403"
$ws.Range("E29").Value = "  -5.19%  "

$ws.Range("D30").Value = "'1.433"
$ws.Range("E30").Value = "  -1.65%  "

$ws.Range("D31").Value = "'3.710"
$ws.Range("E31").Value = "  -3.00%  "

$ws.Range("D32").Value = "'3.684"
$ws.Range("E32").Value = "  -3.88%  "

$ws.Range("D33").Value = "'1.679"
$ws.Range("E33").Value = "  -3.11%  "

$ws.Range("D34").Value = "'1.049"
$ws.Range("E34").Value = "  -6.84%  "

$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").Value = "'0.6299"
$ws.Range("E36").Value = "  -4.10%  "

$ws.Range("D37").Value = "'2.738"
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").Value = "'6.429"
$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("E39").Value = "  -3.85%  "

$ws.Range("D40").Value = "1.130.17"
$ws.Range("E40").Value = "  -8.24%  "

$ws.Range("D41").Value = "'0.8657"
$ws.Range("E41").Value = "  -6.30%  "

$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("D43").Value = "1.972.94"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "'99.56"
$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").Value = "'60.12"
$ws.Range("E45").Value = "  -5.27%  "

$ws.Range("E46").Value = "  -2.51%  "

$ws.Range("D47").Value = "'1.570"
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("D48").Value = "'0.05482"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("D49").Value = "'0.4510"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D50").Value = "'8.231"
$ws.Range("E50").Value = "  -3.38%  "

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.55%  "
